# Borjomi.xlsx upgrade: extend left table with a new "2023/2025" column (P)
# and replace "-" placeholders with "…" ellipsis characters, per commit
# "upgrade left table until javakheti".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Replicate column O's formatting into the new column P for the whole
#    table body (rows 4-18) by copying; we will overwrite values after.
# ---------------------------------------------------------------------
$ws.Range("O4:O18").Copy($ws.Range("P4"))

# Column P should be the same width as the rest of the data columns
# (it used to be a wider "spacer" column).
$ws.Columns("P").ColumnWidth = 10.7109375

# ---------------------------------------------------------------------
# 2) New header cell for column P.
# ---------------------------------------------------------------------
$ws.Range("O3").Copy($ws.Range("P3"))
$ws.Range("P3").Value = "2023/2025"

# ---------------------------------------------------------------------
# 3) Fill in the actual new-column values.
# ---------------------------------------------------------------------
$ws.Range("P4").Value = 20
$ws.Range("P5").Value = 3547

$ws.Range("P6").Value = "…"
$ws.Range("P7").Value = "…"
$ws.Range("P8").Value = "…"
$ws.Range("P9").Value = "…"
$ws.Range("P10").Value = "…"
$ws.Range("P11").Value = "…"
$ws.Range("P12").Value = "…"
$ws.Range("P13").Value = "…"
$ws.Range("P14").Value = "…"
$ws.Range("P15").Value = "…"

$ws.Range("P16").Value = 430
$ws.Range("P17").Value = 376
$ws.Range("P18").Value = 54

# ---------------------------------------------------------------------
# 4) Column O (2023/2024) previously showed "-" for several rows where
#    data is now available; fill in the real figures.
# ---------------------------------------------------------------------
$ws.Range("O6").Value = 299
$ws.Range("O7").Value = 147
$ws.Range("O8").Value = 152
$ws.Range("O9").Value = 273
$ws.Range("O10").Value = 133
$ws.Range("O11").Value = 140
$ws.Range("O12").Value = 261
$ws.Range("O13").Value = 130
$ws.Range("O14").Value = 131
$ws.Range("O15").Value = 11.906931535143672
# This particular figure is shown without the grey band fill, unlike the
# other cells in this column.
$ws.Range("O15").Interior.Pattern = -4142

# ---------------------------------------------------------------------
# 5) Replace every remaining "-" placeholder in the table body with the
#    "…" ellipsis character used for confidential/unavailable data.
# ---------------------------------------------------------------------
$dashRanges = @(
  "B6","B7","B8","B9","B10","B11","B12","B13","B14","B15","B16","B17","B18",
  "C16","D16","E16","F16","G16",
  "C17","D17","E17","F17","G17",
  "C18","D18","E18","F18","G18"
)
foreach ($addr in $dashRanges) {
  $ws.Range($addr).Value = "…"
}

# ---------------------------------------------------------------------
# 6) Add the new footnote row explaining the "…" placeholder.
# ---------------------------------------------------------------------
$ws.Range("A20").Copy($ws.Range("A21"))
$ws.Range("A21").Value = "Note: „ ... „ - Data is confidential or unavailable."
$ws.Range("A21").Font.Name = "Arial"
$ws.Range("A21").Font.Size = 9
$ws.Range("A21").Font.Color = 0
$ws.Range("A21").Font.Bold = $true
$ws.Range("A21").Font.Underline = $true
$note = $ws.Range("A21").Characters(6, 47)
$note.Font.Bold = $false
$note.Font.Underline = $false
